# Update db_rdmarcas sheet: convert row 2 metrics to real numbers and
# append the new sales-tracking rows 3-6 (row 6 keeps its metrics as
# text, matching the original "N.NN" string formatting used upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing row, numeric metrics (was text "5000.00" etc.) ---
$ws.Cells.Item(2, 2).Value = 5000
$ws.Cells.Item(2, 3).Value = 5000
$ws.Cells.Item(2, 4).Value = 5000
$ws.Cells.Item(2, 5).Value = 5000
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 100

# --- Column A for the new rows must stay text (these look like dates
# but are free-form strings in the source data, same as "05/50/5000"
# in A2), so force text formatting before writing them. ---
$ws.Range("A3:A6").NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "08/06/2000"
$ws.Cells.Item(4, 1).Value = "03/08/2023"
$ws.Cells.Item(5, 1).Value = "03/08/2023"
$ws.Cells.Item(6, 1).Value = "03/08/2023"

# --- Row 3: numeric metrics ---
$ws.Cells.Item(3, 2).Value = 4000
$ws.Cells.Item(3, 3).Value = 9000
$ws.Cells.Item(3, 4).Value = 4000
$ws.Cells.Item(3, 5).Value = 9000
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 100

# --- Row 4: numeric metrics ---
$ws.Cells.Item(4, 2).Value = 4000
$ws.Cells.Item(4, 3).Value = 13000
$ws.Cells.Item(4, 4).Value = 5000
$ws.Cells.Item(4, 5).Value = 14000
$ws.Cells.Item(4, 6).Value = 1000
$ws.Cells.Item(4, 7).Value = 107.69

# --- Row 5: numeric metrics ---
$ws.Cells.Item(5, 2).Value = 4740
$ws.Cells.Item(5, 3).Value = 17740
$ws.Cells.Item(5, 4).Value = 4041
$ws.Cells.Item(5, 5).Value = 18041
$ws.Cells.Item(5, 6).Value = 301
$ws.Cells.Item(5, 7).Value = 101.7

# --- Row 6: metrics stored as text (e.g. "4141.00"), so pre-format the
# range as text before writing, otherwise Excel would parse them back
# into numbers. ---
$ws.Range("B6:G6").NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "4141.00"
$ws.Cells.Item(6, 3).Value = "21881.00"
$ws.Cells.Item(6, 4).Value = "4142.00"
$ws.Cells.Item(6, 5).Value = "22183.00"
$ws.Cells.Item(6, 6).Value = "302.00"
$ws.Cells.Item(6, 7).Value = "101.38"
